$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = $null
$ws.Range("H40").Value = 3809.25
$ws.Range("J40").Value = 4191.3335
$ws.Range("L40").Value = 4191.3335
$ws.Range("N40").Value = -4541.3335
$ws.Range("H55").Value = 92.947365
$ws.Range("I55").Value = 41.545456
$ws.Range("J55").Value = 163.625
$ws.Range("K55").Value = 41.545456
$ws.Range("L55").Value = 163.625
$ws.Range("M55").Value = 172.454544
$ws.Range("N55").Value = -591.625
$ws.Range("H103").Value = 993.35297
$ws.Range("J103").Value = 880.7273
$ws.Range("L103").Value = 2642.1819
$ws.Range("N103").Value = -3814.1819
$ws.Range("H112").Value = 1943.3158
$ws.Range("J112").Value = 2148.2307
$ws.Range("L112").Value = 6444.6921
$ws.Range("N112").Value = -8660.6921
$ws.Range("H113").Value = 6482.16
$ws.Range("I113").Value = 6902.7
$ws.Range("J113").Value = 4800
$ws.Range("K113").Value = 6902.7
$ws.Range("L113").Value = 4800
$ws.Range("M113").Value = -3648.7
$ws.Range("N113").Value = -11308

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 2133.1428
$ws.Range("I97").Value = 233
$ws.Range("J97").Value = 4666.6665
$ws.Range("K97").Value = 233
$ws.Range("L97").Value = 4666.6665
$ws.Range("M97").Value = 263
$ws.Range("N97").Value = -5658.6665
$ws.Range("H110").Value = 1354.4445
$ws.Range("J110").Value = 1999
$ws.Range("L110").Value = 1999
$ws.Range("N110").Value = -6089
$ws.Range("H122").Value = 4166.8667
$ws.Range("I122").Value = 4001.1667
$ws.Range("J122").Value = 4829.6665
$ws.Range("K122").Value = 12003.5001
$ws.Range("L122").Value = 14488.9995
$ws.Range("M122").Value = -9553.500100000001
$ws.Range("N122").Value = -19388.9995
$ws.Range("H132").Value = 1789.5
$ws.Range("I132").Value = 1777.0605
$ws.Range("J132").Value = 2200
$ws.Range("K132").Value = 5331.181500000001
$ws.Range("L132").Value = 6600
$ws.Range("M132").Value = -2801.181500000001
$ws.Range("N132").Value = -11660

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1328
$ws.Range("I94").Value = 1565
$ws.Range("J94").Value = 759.2
$ws.Range("K94").Value = 1565
$ws.Range("L94").Value = 759.2
$ws.Range("M94").Value = -1114
$ws.Range("N94").Value = -1661.2
$ws.Range("H99").Value = 3364
$ws.Range("I99").Value = 1705
$ws.Range("K99").Value = 1705
$ws.Range("M99").Value = -207
$ws.Range("H107").Value = 3712.8408
$ws.Range("I107").Value = 1102.3103
$ws.Range("K107").Value = 1102.3103
$ws.Range("M107").Value = 817.6896999999999
$ws.Range("H130").Value = 75000
$ws.Range("J130").Value = 75000
$ws.Range("L130").Value = 75000
$ws.Range("N130").Value = -85040
$ws.Range("H134").Value = 3090.739
$ws.Range("I134").Value = 3157.2942
$ws.Range("K134").Value = 9471.882599999999
$ws.Range("M134").Value = -6936.882599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 361.7143
$ws.Range("J2").Value = 265.6
$ws.Range("L2").Value = 265.6
$ws.Range("N2").Value = -491.6
$ws.Range("H86").Value = 7959.846
$ws.Range("I86").Value = 7565.737
$ws.Range("K86").Value = 7565.737
$ws.Range("M86").Value = -6442.737
$ws.Range("H89").Value = 7959.846
$ws.Range("I89").Value = 7565.737
$ws.Range("K89").Value = 37828.685
$ws.Range("M89").Value = -32212.685
$ws.Range("H99").Value = 3260.3333
$ws.Range("J99").Value = 5500
$ws.Range("L99").Value = 5500
$ws.Range("N99").Value = -8496
$ws.Range("H122").Value = 2743.0527
$ws.Range("I122").Value = 2324.9167
$ws.Range("K122").Value = 6974.750100000001
$ws.Range("M122").Value = -4524.750100000001
$ws.Range("H126").Value = 3260.3333
$ws.Range("J126").Value = 5500
$ws.Range("L126").Value = 16500
$ws.Range("N126").Value = -21440
$ws.Range("H134").Value = 1163.1666
$ws.Range("I134").Value = 1031.4642
$ws.Range("K134").Value = 3094.3926
$ws.Range("M134").Value = -559.3925999999997

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 613.25
$ws.Range("I8").Value = 613.25
$ws.Range("K8").Value = 1839.75
$ws.Range("M8").Value = -1700.75
$ws.Range("H14").Value = 146.08333
$ws.Range("I14").Value = 146.08333
$ws.Range("K14").Value = 438.24999
$ws.Range("M14").Value = -265.24999
$ws.Range("H76").Value = 3509.5
$ws.Range("I76").Value = 3509.5
$ws.Range("K76").Value = 10528.5
$ws.Range("M76").Value = -10145.5
$ws.Range("H79").Value = 3509.5
$ws.Range("I79").Value = 3509.5
$ws.Range("K79").Value = 10528.5
$ws.Range("M79").Value = -9202.5
$ws.Range("H104").Value = 19760
$ws.Range("I104").Value = 9000
$ws.Range("J104").Value = 22450
$ws.Range("K104").Value = 27000
$ws.Range("L104").Value = 67350
$ws.Range("M104").Value = -24379
$ws.Range("N104").Value = -72592
$ws.Range("H118").Value = 257
$ws.Range("I118").Value = 257
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 771
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = 472
$ws.Range("N118").Value = $null

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9256.227999999999
$ws.Range("I70").Value = 11219.6
$ws.Range("K70").Value = 11219.6
$ws.Range("M70").Value = -10949.6
$ws.Range("H73").Value = 9256.227999999999
$ws.Range("I73").Value = 11219.6
$ws.Range("K73").Value = 11219.6
$ws.Range("M73").Value = -10283.6
$ws.Range("H97").Value = 5271.154
$ws.Range("I97").Value = 1602.3182
$ws.Range("J97").Value = 25449.75
$ws.Range("K97").Value = 1602.3182
$ws.Range("L97").Value = 25449.75
$ws.Range("M97").Value = -1106.3182
$ws.Range("N97").Value = -26441.75
$ws.Range("H102").Value = 514.1177
$ws.Range("I102").Value = 578.6667
$ws.Range("K102").Value = 578.6667
$ws.Range("M102").Value = 1043.3333
$ws.Range("H107").Value = 419.45456
$ws.Range("J107").Value = 427.4
$ws.Range("L107").Value = 427.4
$ws.Range("N107").Value = -4267.4
$ws.Range("H132").Value = 2694.6858
$ws.Range("I132").Value = 2268.8667
$ws.Range("J132").Value = 5249.6
$ws.Range("K132").Value = 6806.6001
$ws.Range("L132").Value = 15748.8
$ws.Range("M132").Value = -4276.6001
$ws.Range("N132").Value = -20808.8
$ws.Range("H135").Value = 75000
$ws.Range("J135").Value = 75000
$ws.Range("L135").Value = 75000
$ws.Range("N135").Value = -85140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2358.3333
$ws.Range("I16").Value = 2029.2222
$ws.Range("J16").Value = 4333
$ws.Range("K16").Value = 2029.2222
$ws.Range("L16").Value = 4333
$ws.Range("M16").Value = -1859.2222
$ws.Range("N16").Value = -4673
$ws.Range("H40").Value = 6208.913
$ws.Range("I40").Value = 4089
$ws.Range("K40").Value = 4089
$ws.Range("M40").Value = -3953
$ws.Range("H100").Value = 4150.136
$ws.Range("I100").Value = 2490.4167
$ws.Range("J100").Value = 6141.8
$ws.Range("K100").Value = 2490.4167
$ws.Range("L100").Value = 6141.8
$ws.Range("M100").Value = -1949.4167
$ws.Range("N100").Value = -7223.8
$ws.Range("H132").Value = 5271.6113
$ws.Range("I132").Value = 5383.077
$ws.Range("K132").Value = 16149.231
$ws.Range("M132").Value = -13619.231
$ws.Range("H136").Value = 3645.6785
$ws.Range("I136").Value = 4055.6875
$ws.Range("K136").Value = 12167.0625
$ws.Range("M136").Value = -9617.0625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9100
$ws.Range("I62").Value = 4750
$ws.Range("K62").Value = 4750
$ws.Range("M62").Value = -4126
$ws.Range("H65").Value = 9100
$ws.Range("I65").Value = 4750
$ws.Range("K65").Value = 23750
$ws.Range("M65").Value = -20630
$ws.Range("H107").Value = 1251.2
$ws.Range("I107").Value = 1032.1666
$ws.Range("J107").Value = 1814.4286
$ws.Range("K107").Value = 3096.4998
$ws.Range("L107").Value = 5443.2858
$ws.Range("M107").Value = -1176.4998
$ws.Range("N107").Value = -9283.2858
$ws.Range("H122").Value = 5721.3477
$ws.Range("I122").Value = 3809.3
$ws.Range("K122").Value = 11427.9
$ws.Range("M122").Value = -8977.900000000001
$ws.Range("H126").Value = 3161.875
$ws.Range("I126").Value = 2574.75
$ws.Range("K126").Value = 7724.25
$ws.Range("M126").Value = -5254.25
$ws.Range("H132").Value = 1935.262
$ws.Range("I132").Value = 1832.025
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 5496.075000000001
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -2966.075000000001
$ws.Range("N132").Value = -17060
$ws.Range("H136").Value = 1782.5
$ws.Range("I136").Value = 1615.7407
$ws.Range("K136").Value = 4847.2221
$ws.Range("M136").Value = -2297.2221
